$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.370.35"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.56%  "

$ws.Range("D3").Value = "'1.870.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'243.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.91%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Value = "'0.4709"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.08%  "

$ws.Range("D8").Value = "'0.2876"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.45%  "

$ws.Range("D9").Value = "'0.06452"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("D10").Value = "'21.99"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.08%  "

$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").Value = "'1.868.97"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.27%  "

$ws.Range("D13").Value = "'96.07"
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").Value = "'0.7246"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.90%  "

$ws.Range("D15").Value = "'5.123"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.10%  "

$ws.Range("D16").Value = "'279.33"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.96%  "

$ws.Range("D17").Value = "'30.358.57"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.92%  "

$ws.Range("D18").Value = "'12.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.68%  "

$ws.Range("D19").Value = "'1.000"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.00%  "

$ws.Range("D20").Value = "'0.000007502"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.17%  "

$ws.Range("D21").Value = "'2.116.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.14%  "

$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").Value = "'5.233"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").Value = "'6.225"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.49%  "

$ws.Range("D25").Value = "'163.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.35%  "

$ws.Range("E26").Value = "  -1.55%  "

$ws.Range("D27").Value = "'18.67"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.74%  "

$ws.Range("E28").Value = "  -1.72%  "

$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "'0.09633"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.16%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'1.321"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.41%  "

$ws.Range("D31").Value = "'1.487"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.09%  "

$ws.Range("E32").Value = "  -1.35%  "

$ws.Range("D33").Value = "'4.110"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("D34").Value = "'0.04809"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.10%  "

$ws.Range("D35").Value = "'1.120"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("D36").Value = "'0.6872"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.21%  "

$ws.Range("D37").Value = "'2.706"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.45%  "

$ws.Range("D38").Value = "'0.01879"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.59%  "

$ws.Range("D39").Value = "'2.812"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.85%  "

$ws.Range("E40").Value = "  -0.75%  "

$ws.Range("D41").Value = "'74.09"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.90%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D42").Value = "'1.934"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.36%  "

$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").Value = "'0.4221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("D44").Value = "'0.9999"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.04%  "

$ws.Range("D45").Value = "'0.8257"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.19%  "

$ws.Range("E46").Value = "  -0.87%  "

$ws.Range("D47").Value = "'9.559"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.00%  "

$ws.Range("D48").Value = "'35.26"
$ws.Range("D48").Style = "Normal"

$ws.Range("D49").Value = "'6.953"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.43%  "

$ws.Range("D50").Value = "'901.89"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.21%  "

$ws.Range("D51").Value = "'0.05720"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.78%  "
